$wb = $excel.ActiveWorkbook

# --- Sheet "DatosMotor" (sheet3): update SMP036 -> SMP038 values ---
$wsMotor = $wb.Worksheets.Item("DatosMotor")
$wsMotor.Range("A2").Value = "SMP038"
$wsMotor.Range("B2").Value = "ABC12SSMP038"
$wsMotor.Range("C2").Value = "ZAZ123SSMP038"
$wsMotor.Range("A3:E3").Select()

# --- Sheet "DatosCuenta" (sheet1): update Smoke names ---
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("B2").Value = "SmokeName"
$wsCuenta.Range("A2").Value = "Smoke"
$wsCuenta.Range("A2").Select()

# --- Sheet "DatosHogar" (sheet2): increment numeric value ---
$wsHogar = $wb.Worksheets.Item("DatosHogar")
$wsHogar.Range("A2").Value = 655
$wsHogar.Range("A2").Select()

# --- Sheet "DatosAP" (sheet4): increment numeric value, update selection ---
$wsAP = $wb.Worksheets.Item("DatosAP")
$wsAP.Range("A2").Value = 21200136
$wsAP.Activate()
$wsAP.Range("E10").Select()
